# Updated symbol list on Tue Jan 17 08:56:58 UTC 2023 with GitHub Actions
# Applies updated Price (D) and Volume(1h) (E) values to the crypto sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2;  D = "301.05";     E = "0.49%" },
    @{ Row = 3;  D = "31.46";      E = "0.38%" },
    @{ Row = 4;  D = "5.085";      E = "-0.99%" },
    @{ Row = 5;  D = "0.07854";    E = "-5.10%" },
    @{ Row = 6;  D = "2.350";      E = "0.39%" },
    @{ Row = 7;  D = "7.806";      E = "-1.24%" },
    @{ Row = 8;  D = "3.820";      E = "-0.74%" },
    @{ Row = 9;  D = "0.9213" },
    @{ Row = 10; D = "0.1755";     E = "1.73%" },
    @{ Row = 11; E = "2.94%" },
    @{ Row = 12; D = "0.09221";    E = "15.07%" },
    @{ Row = 13; D = "0.02995";    E = "-0.98%" },
    @{ Row = 14; D = "0.1002";     E = "0.68%" },
    @{ Row = 15; D = "0.001511";   E = "-1.00%" },
    @{ Row = 16; D = "0.005815";   E = "-5.27%" },
    @{ Row = 17; D = "3.474";      E = "-0.80%" },
    @{ Row = 18; E = "-0.33%" },
    @{ Row = 20; D = "0.1288";     E = "-4.61%" },
    @{ Row = 21; D = "4.085";      E = "-11.34%" },
    @{ Row = 22; E = "6.36%" },
    @{ Row = 23; D = "0.04600";    E = "-0.07%" },
    @{ Row = 24; D = "0.001249";   E = "-1.45%" },
    @{ Row = 25; D = "0.004472";   E = "0.44%" },
    @{ Row = 26; E = "5.27%" },
    @{ Row = 27; E = "-1.80%" },
    @{ Row = 39; D = "0.01756";    E = "-3.90%" },
    @{ Row = 40; D = "0.04695";    E = "3.76%" },
    @{ Row = 41; D = "0.007114";   E = "-2.03%" },
    @{ Row = 42; D = "0.1358";     E = "0.99%" },
    @{ Row = 43; E = "-0.20%" },
    @{ Row = 44; D = "0.009760";   E = "-8.03%" },
    @{ Row = 45; D = "0.00006260"; E = "-0.48%" },
    @{ Row = 46; E = "-0.62%" },
    @{ Row = 47; E = "19.72%" },
    @{ Row = 48; D = "0.7439";     E = "-9.34%" },
    @{ Row = 49; D = "0.00002096"; E = "-0.62%" },
    @{ Row = 50; D = "0.0001996";  E = "-0.62%" }
)

foreach ($u in $updates) {
    if ($u.ContainsKey("D")) {
        $cell = $ws.Range("D$($u.Row)")
        $cell.NumberFormat = "@"
        $cell.Value = $u.D
        $cell.Style = "Normal"
    }
    if ($u.ContainsKey("E")) {
        $cell = $ws.Range("E$($u.Row)")
        $cell.NumberFormat = "@"
        $cell.Value = $u.E
        $cell.Style = "Normal"
    }
}
